# Update the "Förändrad" (Changed) date column (C) for rows 2-28
# from 45204 (2023-10-05) to 45207 (2023-10-08), keeping all other
# cell data/formatting untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value2 = 45207
    }
}
